$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.051757249440491
$ws.Range("D2").Value2 = 1.061732250817522
$ws.Range("E2").Value2 = 1.058865471869162
$ws.Range("F2").Value2 = 1.070310756288933
$ws.Range("I2").Value2 = 1.02359499962809
$ws.Range("J2").Value2 = 1.056783439289323
$ws.Range("K2").Value2 = 1.064455512831643
$ws.Range("L2").Value2 = 1.061596538619531
$ws.Range("M2").Value2 = 1.073010933640486
$ws.Range("N2").Value2 = 1.058284192987868
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.053346629947035
$ws.Range("D3").Value2 = 1.06326392973464
$ws.Range("E3").Value2 = 1.060309036370354
$ws.Range("F3").Value2 = 1.071903066491194
$ws.Range("I3").Value2 = 1.023504579208684
$ws.Range("J3").Value2 = 1.058019915742132
$ws.Range("K3").Value2 = 1.065800065445098
$ws.Range("L3").Value2 = 1.062852626451615
$ws.Range("M3").Value2 = 1.074417661830046
$ws.Range("N3").Value2 = 1.059522425379066
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.05437422483451
$ws.Range("D4").Value2 = 1.064254536129161
$ws.Range("E4").Value2 = 1.061242541667407
$ws.Range("F4").Value2 = 1.072933029893937
$ws.Range("I4").Value2 = 1.023443335729026
$ws.Range("J4").Value2 = 1.058818792020995
$ws.Range("K4").Value2 = 1.06666908702309
$ws.Range("L4").Value2 = 1.063664303075816
$ws.Range("M4").Value2 = 1.075327047972722
$ws.Range("N4").Value2 = 1.060322436153878
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.054806032897938
$ws.Range("D5").Value2 = 1.06467087651156
$ws.Range("E5").Value2 = 1.061634855790218
$ws.Range("F5").Value2 = 1.073365946019978
$ws.Range("I5").Value2 = 1.023416932628374
$ws.Range("J5").Value2 = 1.059154356991065
$ws.Range("K5").Value2 = 1.067034192431198
$ws.Range("L5").Value2 = 1.064005275694892
$ws.Range("M5").Value2 = 1.075709154653179
$ws.Range("N5").Value2 = 1.060658477664695
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.054878524308068
$ws.Range("D6").Value2 = 1.064740775623078
$ws.Range("E6").Value2 = 1.061700719541604
$ws.Range("F6").Value2 = 1.073438630041589
$ws.Range("I6").Value2 = 1.023412460900047
$ws.Range("J6").Value2 = 1.059210683434686
$ws.Range("K6").Value2 = 1.06709548186769
$ws.Range("L6").Value2 = 1.0640625116519
$ws.Range("M6").Value2 = 1.075773300561373
$ws.Range("N6").Value2 = 1.060714884098327
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.054379995422263
$ws.Range("D7").Value2 = 1.064260099715058
$ws.Range("E7").Value2 = 1.061247784296302
$ws.Range("F7").Value2 = 1.072938814847315
$ws.Range("I7").Value2 = 1.023442985510381
$ws.Range("J7").Value2 = 1.058823276959313
$ws.Range("K7").Value2 = 1.066673966479887
$ws.Range("L7").Value2 = 1.06366886016781
$ws.Range("M7").Value2 = 1.075332154475309
$ws.Range("N7").Value2 = 1.060326927461322
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.052294567142697
$ws.Range("D8").Value2 = 1.062249995217486
$ws.Range("E8").Value2 = 1.0593534548851
$ws.Range("F8").Value2 = 1.07084896568446
$ws.Range("I8").Value2 = 1.023565007130014
$ws.Range("J8").Value2 = 1.057201565920324
$ws.Range("K8").Value2 = 1.064910120493698
$ws.Range("L8").Value2 = 1.062021270095994
$ws.Range("M8").Value2 = 1.073486526459869
$ws.Range("N8").Value2 = 1.058702913406646
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.048612961053824
$ws.Range("D9").Value2 = 1.058703836005057
$ws.Range("E9").Value2 = 1.056010679151949
$ws.Range("F9").Value2 = 1.067163228654945
$ws.Range("I9").Value2 = 1.023759139049229
$ws.Range("J9").Value2 = 1.054334381240272
$ws.Range("K9").Value2 = 1.061794087246479
$ws.Range("L9").Value2 = 1.059109332154391
$ws.Range("M9").Value2 = 1.070227390372366
$ws.Range("N9").Value2 = 1.055831656995477
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.046153446737528
$ws.Range("D10").Value2 = 1.056336514650688
$ws.Range("E10").Value2 = 1.053778560551304
$ws.Range("F10").Value2 = 1.064703473515629
$ws.Range("I10").Value2 = 1.023874603872166
$ws.Range("J10").Value2 = 1.052416128232756
$ws.Range("K10").Value2 = 1.059710992844307
$ws.Range("L10").Value2 = 1.057161818270381
$ws.Range("M10").Value2 = 1.068049558024812
$ws.Range("N10").Value2 = 1.053910679848663
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.045087120018198
$ws.Range("D11").Value2 = 1.055310573094481
$ws.Range("E11").Value2 = 1.052811079498881
$ws.Range("F11").Value2 = 1.063637650253416
$ws.Range("I11").Value2 = 1.023921309042974
$ws.Range("J11").Value2 = 1.051583803408769
$ws.Range("K11").Value2 = 1.058807535004636
$ws.Range("L11").Value2 = 1.056316960302155
$ws.Range("M11").Value2 = 1.067105232086856
$ws.Range("N11").Value2 = 1.053077173027958
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.04469082716552
$ws.Range("D12").Value2 = 1.05492935224563
$ws.Range("E12").Value2 = 1.052451561657781
$ws.Range("F12").Value2 = 1.06324163671049
$ws.Range("I12").Value2 = 1.023938164268129
$ws.Range("J12").Value2 = 1.051274376567358
$ws.Range("K12").Value2 = 1.058471722419086
$ws.Range("L12").Value2 = 1.056002898748741
$ws.Range("M12").Value2 = 1.066754262157248
$ws.Range("N12").Value2 = 1.05276730676494
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.044775843068781
$ws.Range("D13").Value2 = 1.055011131912082
$ws.Range("E13").Value2 = 1.052528686479753
$ws.Range("F13").Value2 = 1.063326588507862
$ws.Range("I13").Value2 = 1.023934571058996
$ws.Range("J13").Value2 = 1.051340761760534
$ws.Range("K13").Value2 = 1.058543765812633
$ws.Range("L13").Value2 = 1.056070277167501
$ws.Range("M13").Value2 = 1.066829555821012
$ws.Range("N13").Value2 = 1.052833786232705
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.045054366703779
$ws.Range("D14").Value2 = 1.055279064177523
$ws.Range("E14").Value2 = 1.052781364796673
$ws.Range("F14").Value2 = 1.063604918138076
$ws.Range("I14").Value2 = 1.023922712353276
$ws.Range("J14").Value2 = 1.051558231526436
$ws.Range("K14").Value2 = 1.058779781332342
$ws.Range("L14").Value2 = 1.056291004891293
$ws.Range("M14").Value2 = 1.067076225038378
$ws.Range("N14").Value2 = 1.053051564830619
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.045225945958518
$ws.Range("D15").Value2 = 1.055444127207645
$ws.Range("E15").Value2 = 1.052937027818698
$ws.Range("F15").Value2 = 1.063776390121102
$ws.Range("I15").Value2 = 1.023915340507297
$ws.Range("J15").Value2 = 1.051692186569136
$ws.Range("K15").Value2 = 1.058925167811431
$ws.Range("L15").Value2 = 1.056426970008624
$ws.Range("M15").Value2 = 1.06722817863725
$ws.Range("N15").Value2 = 1.053185710104844
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.046224186212693
$ws.Range("D16").Value2 = 1.056404583803001
$ws.Range("E16").Value2 = 1.053842748177147
$ws.Range("F16").Value2 = 1.064774192467354
$ws.Range("I16").Value2 = 1.023871434968932
$ws.Range("J16").Value2 = 1.052471330240455
$ws.Range("K16").Value2 = 1.059770920823153
$ws.Range("L16").Value2 = 1.057217854954323
$ws.Range("M16").Value2 = 1.068112201508568
$ws.Range("N16").Value2 = 1.053965960249545
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.046849989808693
$ws.Range("D17").Value2 = 1.057006812114231
$ws.Range("E17").Value2 = 1.054410620095456
$ws.Range("F17").Value2 = 1.065399884786179
$ws.Range("I17").Value2 = 1.023843013862359
$ws.Range("J17").Value2 = 1.052959603656944
$ws.Range("K17").Value2 = 1.060301041604037
$ws.Range("L17").Value2 = 1.057713530093707
$ws.Range("M17").Value2 = 1.068666369025145
$ws.Range("N17").Value2 = 1.054454927070292
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.047214881783904
$ws.Range("D18").Value2 = 1.057357997803006
$ws.Range("E18").Value2 = 1.054741758210292
$ws.Range("F18").Value2 = 1.065764770425395
$ws.Range("I18").Value2 = 1.023826118601784
$ws.Range("J18").Value2 = 1.053244240942755
$ws.Range("K18").Value2 = 1.06061011136363
$ws.Range("L18").Value2 = 1.05800249789374
$ws.Range("M18").Value2 = 1.068989479720653
$ws.Range("N18").Value2 = 1.054739968573697
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.047339278999492
$ws.Range("D19").Value2 = 1.057477729085189
$ws.Range("E19").Value2 = 1.054854652398585
$ws.Range("F19").Value2 = 1.065889175384912
$ws.Range("I19").Value2 = 1.023820303822167
$ws.Range("J19").Value2 = 1.053341267257324
$ws.Range("K19").Value2 = 1.060715472540577
$ws.Range("L19").Value2 = 1.058101003040974
$ws.Range("M19").Value2 = 1.069099631025886
$ws.Range("N19").Value2 = 1.054837132676761
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.04678286038866
$ws.Range("D20").Value2 = 1.056942207458154
$ws.Range("E20").Value2 = 1.054349702385182
$ws.Range("F20").Value2 = 1.065332761239748
$ws.Range("I20").Value2 = 1.023846096030143
$ws.Range("J20").Value2 = 1.052907233593889
$ws.Range("K20").Value2 = 1.060244179275873
$ws.Range("L20").Value2 = 1.057660364528482
$ws.Range("M20").Value2 = 1.068606925172391
$ws.Range("N20").Value2 = 1.05440248263574
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.044972354359534
$ws.Range("D21").Value2 = 1.05520016875294
$ws.Range("E21").Value2 = 1.052706961600793
$ws.Range("F21").Value2 = 1.06352296036734
$ws.Range("I21").Value2 = 1.02392621804881
$ws.Range("J21").Value2 = 1.051494199459455
$ws.Range("K21").Value2 = 1.058710287007089
$ws.Range("L21").Value2 = 1.056226012866137
$ws.Range("M21").Value2 = 1.067003592834982
$ws.Range("N21").Value2 = 1.052987441830759
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.043832787831138
$ws.Range("D22").Value2 = 1.054104062682096
$ws.Range("E22").Value2 = 1.051673219985603
$ws.Range("F22").Value2 = 1.062384371735761
$ws.Range("I22").Value2 = 1.023973741382756
$ws.Range("J22").Value2 = 1.050604235588034
$ws.Range("K22").Value2 = 1.057744544170503
$ws.Range("L22").Value2 = 1.055322765362694
$ws.Range("M22").Value2 = 1.065994322560932
$ws.Range("N22").Value2 = 1.052096214108559
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.044437013574077
$ws.Range("D23").Value2 = 1.054685209695664
$ws.Range("E23").Value2 = 1.05222131265032
$ws.Range("F23").Value2 = 1.062988028074642
$ws.Range("I23").Value2 = 1.023948818288664
$ws.Range("J23").Value2 = 1.051076170171822
$ws.Range("K23").Value2 = 1.058256630939724
$ws.Range("L23").Value2 = 1.055801730495723
$ws.Range("M23").Value2 = 1.066529471595368
$ws.Range("N23").Value2 = 1.052568818893588
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.046813193681789
$ws.Range("D24").Value2 = 1.056971399780816
$ws.Range("E24").Value2 = 1.054377228760683
$ws.Range("F24").Value2 = 1.065363091697928
$ws.Range("I24").Value2 = 1.023844704313072
$ws.Range("J24").Value2 = 1.052930897879111
$ws.Range("K24").Value2 = 1.06026987334938
$ws.Range("L24").Value2 = 1.057684388226073
$ws.Range("M24").Value2 = 1.068633785679766
$ws.Range("N24").Value2 = 1.054426180526962
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.049565606285848
$ws.Range("D25").Value2 = 1.059621135994138
$ws.Range("E25").Value2 = 1.056875470313182
$ws.Range("F25").Value2 = 1.068116502766108
$ws.Range("I25").Value2 = 1.02371141742819
$ws.Range("J25").Value2 = 1.055076786112191
$ws.Range("K25").Value2 = 1.059863205819264
$ws.Range("L25").Value2 = 1.059522425379066
$ws.Range("M25").Value2 = 1.071070817747026
$ws.Range("N25").Value2 = 1.056575116167469
